$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to row 458 (row 1 is the header).
# Column C holds an "updated" date serial value that increments by one
# day (45177 -> 45178) for every data row.
$lastRow = 458

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
